$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 350.06668
$ws.Range("J2").Value = 313.4
$ws.Range("L2").Value = 313.4
$ws.Range("N2").Value = -539.4

$ws.Range("H10").Value = 9582.666999999999
$ws.Range("I10").Value = 8750
$ws.Range("J10").Value = 9999
$ws.Range("K10").Value = 8750
$ws.Range("L10").Value = 9999
$ws.Range("M10").Value = -8457
$ws.Range("N10").Value = -10585

$ws.Range("H39").Value = 285.75
$ws.Range("I39").Value = 252.5
$ws.Range("J39").Value = 385.5
$ws.Range("K39").Value = 757.5
$ws.Range("L39").Value = 1156.5
$ws.Range("M39").Value = -461.5
$ws.Range("N39").Value = -1748.5

$ws.Range("H106").Value = 1220.2727
$ws.Range("I106").Value = 952.3
$ws.Range("J106").Value = 3900
$ws.Range("K106").Value = 952.3
$ws.Range("L106").Value = 3900
$ws.Range("M106").Value = -321.3
$ws.Range("N106").Value = -5162

$ws.Range("H107").Value = 5201
$ws.Range("I107").Value = 4666.3335
$ws.Range("J107").Value = 6003
$ws.Range("K107").Value = 4666.3335
$ws.Range("L107").Value = 6003
$ws.Range("M107").Value = -2746.3335
$ws.Range("N107").Value = -9843

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H112").Value = 1737
$ws.Range("I112").Value = 1133
$ws.Range("J112").Value = 1995.8572
$ws.Range("K112").Value = 3399
$ws.Range("L112").Value = 5987.571599999999
$ws.Range("M112").Value = -2291
$ws.Range("N112").Value = -8203.571599999999

$ws.Range("H137").Value = 1581.5
$ws.Range("I137").Value = 1384.8235
$ws.Range("J137").Value = 2250.2
$ws.Range("K137").Value = 4154.470499999999
$ws.Range("L137").Value = 6750.599999999999
$ws.Range("M137").Value = -1604.470499999999
$ws.Range("N137").Value = -11850.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 100000
$ws.Range("J3").Value = 100000
$ws.Range("L3").Value = 100000
$ws.Range("N3").Value = -100230

$ws.Range("H32").Value = 2780.5193
$ws.Range("I32").Value = 2472.1836
$ws.Range("K32").Value = 2472.1836
$ws.Range("M32").Value = -2185.1836

$ws.Range("H61").Value = 3592.1428
$ws.Range("I61").Value = 3599.2307
$ws.Range("K61").Value = 3599.2307
$ws.Range("M61").Value = -3387.2307

$ws.Range("H74").Value = 9382.154
$ws.Range("I74").Value = 1827.625
$ws.Range("K74").Value = 1827.625
$ws.Range("M74").Value = -953.625

$ws.Range("H77").Value = 9382.154
$ws.Range("I77").Value = 1827.625
$ws.Range("K77").Value = 9138.125
$ws.Range("M77").Value = -4770.125

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H136").Value = 3592.1428
$ws.Range("I136").Value = 3599.2307
$ws.Range("K136").Value = 10797.6921
$ws.Range("M136").Value = -8247.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9546.083000000001
$ws.Range("I20").Value = 7796.4287
$ws.Range("J20").Value = 11995.6
$ws.Range("K20").Value = 7796.4287
$ws.Range("L20").Value = 11995.6
$ws.Range("M20").Value = -7549.4287
$ws.Range("N20").Value = -12489.6

$ws.Range("H31").Value = 126
$ws.Range("J31").Value = 126
$ws.Range("L31").Value = 126
$ws.Range("N31").Value = -630

$ws.Range("H82").Value = 28347.25
$ws.Range("I82").Value = 11129.667
$ws.Range("J82").Value = 80000
$ws.Range("K82").Value = 11129.667
$ws.Range("L82").Value = 80000
$ws.Range("M82").Value = -10746.667
$ws.Range("N82").Value = -80766

$ws.Range("H85").Value = 28347.25
$ws.Range("I85").Value = 11129.667
$ws.Range("J85").Value = 80000
$ws.Range("K85").Value = 11129.667
$ws.Range("L85").Value = 80000
$ws.Range("M85").Value = -9803.666999999999
$ws.Range("N85").Value = -82652

$ws.Range("H94").Value = 1819.0435
$ws.Range("J94").Value = 2153.875
$ws.Range("L94").Value = 2153.875
$ws.Range("N94").Value = -3055.875

$ws.Range("H95").Value = 30600
$ws.Range("J95").Value = 30600
$ws.Range("L95").Value = 30600
$ws.Range("N95").Value = -36092

$ws.Range("H106").Value = 59441.5
$ws.Range("J106").Value = 59441.5
$ws.Range("L106").Value = 59441.5
$ws.Range("N106").Value = -61965.5

$ws.Range("H107").Value = 6156.7144
$ws.Range("J107").Value = 6766.3335
$ws.Range("L107").Value = 6766.3335
$ws.Range("N107").Value = -10606.3335

$ws.Range("H134").Value = 2485.762
$ws.Range("I134").Value = 2120
$ws.Range("K134").Value = 6360
$ws.Range("M134").Value = -3825

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 4750
$ws.Range("I45").Value = 3500
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 3500
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -2907
$ws.Range("N45").Value = -7186

$ws.Range("H105").Value = 1884.3684
$ws.Range("J105").Value = 2005.5714
$ws.Range("L105").Value = 2005.5714
$ws.Range("N105").Value = -5499.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.153847
$ws.Range("I2").Value = 28.470589
$ws.Range("K2").Value = 170.823534
$ws.Range("M2").Value = -57.823534

$ws.Range("H4").Value = 3233000
$ws.Range("I4").Value = 3392100
$ws.Range("J4").Value = 51000
$ws.Range("K4").Value = 10176300
$ws.Range("L4").Value = 153000
$ws.Range("M4").Value = -10176188
$ws.Range("N4").Value = -153224

$ws.Range("H11").Value = 66717.87
$ws.Range("J11").Value = 200075.8
$ws.Range("L11").Value = 600227.3999999999
$ws.Range("N11").Value = -600507.3999999999

$ws.Range("H23").Value = 190.4375
$ws.Range("I23").Value = 75.61539
$ws.Range("J23").Value = 688
$ws.Range("K23").Value = 226.84617
$ws.Range("L23").Value = 2064
$ws.Range("M23").Value = 8.153829999999971
$ws.Range("N23").Value = -2534

$ws.Range("H26").Value = 1593.125
$ws.Range("I26").Value = 153.63637
$ws.Range("J26").Value = 4760
$ws.Range("K26").Value = 460.90911
$ws.Range("L26").Value = 14280
$ws.Range("M26").Value = -172.90911
$ws.Range("N26").Value = -14856

$ws.Range("H44").Value = 397
$ws.Range("I44").Value = 397
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1191
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -793
$ws.Range("N44").ClearContents()

$ws.Range("H68").Value = 895
$ws.Range("I68").Value = 895
$ws.Range("K68").Value = 2685
$ws.Range("M68").Value = -1874

$ws.Range("H71").Value = 895
$ws.Range("I71").Value = 895
$ws.Range("K71").Value = 8055
$ws.Range("M71").Value = -3999

$ws.Range("H98").Value = 266
$ws.Range("I98").Value = 111
$ws.Range("K98").Value = 333
$ws.Range("M98").Value = 1165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 100000
$ws.Range("I5").Value = 100000
$ws.Range("K5").Value = 100000
$ws.Range("M5").Value = -99888

$ws.Range("H74").Value = 40001
$ws.Range("J74").Value = 40001
$ws.Range("L74").Value = 40001
$ws.Range("N74").Value = -41873

$ws.Range("H77").Value = 40001
$ws.Range("J77").Value = 40001
$ws.Range("L77").Value = 120003
$ws.Range("N77").Value = -129363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2799998.5
$ws.Range("J2").Value = 2799998.5
$ws.Range("L2").Value = 2799998.5
$ws.Range("N2").Value = -2800222.5

$ws.Range("H22").Value = 2200.625
$ws.Range("I22").Value = 1903
$ws.Range("J22").Value = 2432.111
$ws.Range("K22").Value = 1903
$ws.Range("L22").Value = 2432.111
$ws.Range("M22").Value = -1608
$ws.Range("N22").Value = -3022.111

$ws.Range("H27").Value = 2200.625
$ws.Range("I27").Value = 1903
$ws.Range("J27").Value = 2432.111
$ws.Range("K27").Value = 1903
$ws.Range("L27").Value = 2432.111
$ws.Range("M27").Value = -1796
$ws.Range("N27").Value = -2646.111

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H132").Value = 5212.7896
$ws.Range("I132").Value = 4999.3335
$ws.Range("J132").Value = 5578.7144
$ws.Range("K132").Value = 14998.0005
$ws.Range("L132").Value = 16736.1432
$ws.Range("M132").Value = -12468.0005
$ws.Range("N132").Value = -21796.1432

$ws.Range("H136").Value = 7312.875
$ws.Range("I136").Value = 7312.875
$ws.Range("K136").Value = 21938.625
$ws.Range("M136").Value = -19388.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25249.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H62").Value = 20304.666
$ws.Range("I62").Value = 19749.973
$ws.Range("K62").Value = 19749.973
$ws.Range("M62").Value = -19125.973

$ws.Range("H65").Value = 20304.666
$ws.Range("I65").Value = 19749.973
$ws.Range("K65").Value = 98749.86500000001
$ws.Range("M65").Value = -95629.86500000001

$ws.Range("H132").Value = 2270.3794
$ws.Range("I132").Value = 2113.76
$ws.Range("K132").Value = 6341.280000000001
$ws.Range("M132").Value = -3811.280000000001

$ws.Range("H136").Value = 2891.3845
$ws.Range("I136").Value = 2949.5
$ws.Range("K136").Value = 8848.5
$ws.Range("M136").Value = -6298.5
